# Sample data: add category with dataset
#
# The source CSV-derived workbook had generic header names
# (commune_code / school_code / school_name_en / school_name_km).
# Rename them to the new, more generic category/dataset header names
# (location_id / code / name_en / name_km) on both sheets, and update
# the active sheet / selection to reflect where the author left off
# editing ("Koh Kong" sheet, cell D9 selected).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "Koh Kong"
$ws2 = $wb.Worksheets.Item(2)   # "Kratie"

# Update header row (row 1) on both sheets with the new column names.
$ws1.Range("A1").Value = "location_id"
$ws1.Range("B1").Value = "code"
$ws1.Range("C1").Value = "name_en"
$ws1.Range("D1").Value = "name_km"

$ws2.Range("A1").Value = "location_id"
$ws2.Range("B1").Value = "code"
$ws2.Range("C1").Value = "name_en"
$ws2.Range("D1").Value = "name_km"

# Make "Koh Kong" the active sheet/tab, with D9 as the selected cell.
[void]$ws1.Activate()
[void]$ws1.Range("D9").Select()
